# Updated GSC export data: append two more daily rows (2025-11-13 and
# 2025-11-14) to the "Chart" sheet's breadcrumb table.
#
# Chart sheet layout: column A = Date, column B = Invalid count,
# column C = Valid count. The table currently ends at row 39
# (2025-11-12); we extend it with rows 40 and 41.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$newRows = @(
    @{ Row = 40; Date = "2025-11-13"; Invalid = 0.0; Valid = 41.0 },
    @{ Row = 41; Date = "2025-11-14"; Invalid = 0.0; Valid = 38.0 }
)

foreach ($r in $newRows) {
    $dateCell = $ws.Cells.Item($r.Row, 1)

    # Force the date string to be stored as text (matching the rest of
    # the column) instead of letting Excel auto-convert it to a date
    # serial number, then drop the temporary format override so the
    # cell keeps the sheet's default style.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r.Date
    $dateCell.ClearFormats()

    $ws.Cells.Item($r.Row, 2).Value = $r.Invalid
    $ws.Cells.Item($r.Row, 3).Value = $r.Valid
}
